# Atualização de bases das ligas, do dia: 17-02-2024 às 11:11
#
# Two pairs of match rows had their data swapped between rows (the two
# fixtures played on the same date traded places): row 48 <-> row 49, and
# row 83 <-> row 84. Column A (the positional id) and column E (Date) stay
# put; every other column (matchid, HomeTeam, AwayTeam, score, result, and
# all odds columns) moves with its match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 48 / 49 swap -----------------------------------------------
# Row 48 becomes what row 49 used to be:
$ws.Range("B48").Value = 7035048
$ws.Range("F48").Value = "SG Unterrath"
$ws.Range("G48").Value = "TuRU Dsseldorf"
$ws.Range("H48").Value = 1
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = "H"
$ws.Range("K48").Value = 3.25
$ws.Range("L48").Value = 4
$ws.Range("M48").Value = 1.8
$ws.Range("N48").Value = 2.9
$ws.Range("O48").Value = 4
$ws.Range("P48").Value = 1.95
$ws.Range("Q48").Value = 0.5
$ws.Range("R48").Value = 1.8
$ws.Range("S48").Value = 2
$ws.Range("T48").Value = 3
$ws.Range("U48").Value = 1.75
$ws.Range("V48").Value = 1.95
$ws.Range("W48").Value = 1.9
$ws.Range("X48").Value = -1
$ws.Range("Y48").Value = -1
$ws.Range("Z48").Value = 0.8
$ws.Range("AA48").Value = -1
$ws.Range("AB48").Value = -1
$ws.Range("AC48").Value = 0.95

# Row 49 becomes what row 48 used to be:
$ws.Range("B49").Value = 7035047
$ws.Range("F49").Value = "SC Dsseldorf West"
$ws.Range("G49").Value = "VfL Viktoria JuchenGarzweiler"
$ws.Range("H49").Value = 3
$ws.Range("I49").Value = 4
$ws.Range("J49").Value = "A"
$ws.Range("K49").Value = 1.909
$ws.Range("L49").Value = 3.75
$ws.Range("M49").Value = 3.1
$ws.Range("N49").Value = 2.2
$ws.Range("O49").Value = 3.6
$ws.Range("P49").Value = 2.625
$ws.Range("Q49").Value = -0.25
$ws.Range("R49").Value = 2
$ws.Range("S49").Value = 1.8
$ws.Range("T49").Value = 3
$ws.Range("U49").Value = 1.825
$ws.Range("V49").Value = 1.975
$ws.Range("W49").Value = -1
$ws.Range("X49").Value = -1
$ws.Range("Y49").Value = 1.625
$ws.Range("Z49").Value = -1
$ws.Range("AA49").Value = 0.8
$ws.Range("AB49").Value = 0.825
$ws.Range("AC49").Value = -1

# --- Rows 83 / 84 swap -------------------------------------------------
# Row 83 becomes what row 84 used to be:
$ws.Range("B83").Value = 7511940
$ws.Range("F83").Value = "ASV Suchteln"
$ws.Range("G83").Value = "Holzheimer SG"
$ws.Range("H83").Value = 3
$ws.Range("I83").Value = 3
$ws.Range("J83").Value = "D"
$ws.Range("K83").Value = 2.75
$ws.Range("L83").Value = 3.6
$ws.Range("M83").Value = 2.1
$ws.Range("N83").Value = 4
$ws.Range("O83").Value = 4
$ws.Range("P83").Value = 1.6
$ws.Range("Q83").Value = 1
$ws.Range("R83").Value = 1.8
$ws.Range("S83").Value = 2
$ws.Range("T83").Value = 3.5
$ws.Range("U83").Value = 1.85
$ws.Range("V83").Value = 1.95
$ws.Range("W83").Value = -1
$ws.Range("X83").Value = 3
$ws.Range("Y83").Value = -1
$ws.Range("Z83").Value = 0.8
$ws.Range("AA83").Value = -1
$ws.Range("AB83").Value = 0.8500000000000001
$ws.Range("AC83").Value = -1

# Row 84 becomes what row 83 used to be:
$ws.Range("B84").Value = 7511941
$ws.Range("F84").Value = "SpVgg SterkradeNord"
$ws.Range("G84").Value = "BlauWeiss Mintard"
$ws.Range("H84").Value = 3
$ws.Range("I84").Value = 1
$ws.Range("J84").Value = "H"
$ws.Range("K84").Value = 3.4
$ws.Range("L84").Value = 4.2
$ws.Range("M84").Value = 1.727
$ws.Range("N84").Value = 2.8
$ws.Range("O84").Value = 4.2
$ws.Range("P84").Value = 1.95
$ws.Range("Q84").Value = 0.5
$ws.Range("R84").Value = 1.8
$ws.Range("S84").Value = 2
$ws.Range("T84").Value = 3.5
$ws.Range("U84").Value = 1.875
$ws.Range("V84").Value = 1.925
$ws.Range("W84").Value = 1.8
$ws.Range("X84").Value = -1
$ws.Range("Y84").Value = -1
$ws.Range("Z84").Value = 0.8
$ws.Range("AA84").Value = -1
$ws.Range("AB84").Value = 0.875
$ws.Range("AC84").Value = -1
